# Update simulation result values on Sheet1 as part of the analysis
# of the first 3 simulations ("analisis de las 3 primeras simulaciones").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "I5"  = 95.91309585481312
    "I6"  = 95.93606495153165
    "I7"  = 95.93396177026902
    "G20" = 97.72876282784233
    "G21" = 97.81628922157634
    "G22" = 97.77393150350308
    "H23" = 97.30063949819609
    "H24" = 97.30207358240769
    "H25" = 97.29438423966143
    "I28" = 95.83124397915222
    "I29" = 95.86800017832331
    "G38" = 97.72911929903448
    "G39" = 97.73341690771791
    "H40" = 97.31303137772166
    "H41" = 97.29695221162127
    "I44" = 95.86728050753651
    "I45" = 95.93715225614618
    "G54" = 97.72106763984402
    "G55" = 97.76249721655044
    "H56" = 97.2521272267546
    "H57" = 97.17059767823359
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
